$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Playlist names: append "1" to both GloriaJeans and Starbucks values.
$ws.Range("B5").Value = "GloriaJeans1"
$ws.Range("B6").Value = "Starbucks1"

# Move the active selection to C6 (matches the saved selection in the sheet view).
$ws.Range("C6").Select()
